# Apply "Initial Wk 2 material" edit to the schedule worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Week 1): topic text changes; drop the assignment entry in J2.
$ws.Range("F2").Value = "The Idea of Your Computer: The file system; the terminal; the Unix way of thinking"
$ws.Range("J2").ClearContents()

# Row 3 (Week 2): new topic text, plus new content/example links.
$ws.Range("F3").Value = "The Idea of a Shell: Finding, listing, and inspecting things"
$ws.Range("H3").Value = "/content/02-content"
$ws.Range("I3").Value = "/example/02-example"

# Rows 4-7 keep the same topic text (only the shared-string ordering changed upstream,
# which is not something we set directly via the object model).
$ws.Range("F4").Value = "The Idea of Plain Text: Text editors; slicing and dicing; regular expressions"
$ws.Range("F5").Value = "The Idea of Version Control: Git and GitHub; knowing what you did"
$ws.Range("F6").Value = "The Idea of a Build System: Make, targets, and Quarto; IDEs"
$ws.Range("F7").Value = "The Idea of the Network: Servers, websites, and APIs"

# Widen column F to fit the longer text, and move the active selection to F3.
$ws.Columns.Item(6).ColumnWidth = 71
$ws.Range("F3").Select()
